# Update the cryptocurrency Price (column D) and Volume(1h) (column E) values
# reported on this worksheet to the refreshed figures from the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.820.56'
$ws.Range("E2").Value = '  +1.79%  '
$ws.Range("D3").Value = '3.117.15'
$ws.Range("E3").Value = '  +1.51%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''533.16'
$ws.Range("E5").Value = '  +2.65%  '
$ws.Range("D6").Value = '''138.47'
$ws.Range("E6").Value = '  +2.04%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("E8").Value = '  +10.35%  '
$ws.Range("D9").Value = '''7.35'
$ws.Range("E9").Value = '  +0.38%  '
$ws.Range("E10").Value = '  +1.69%  '
$ws.Range("D11").Value = '''0.413'
$ws.Range("E11").Value = '  +4.63%  '
$ws.Range("D12").Value = '''0.139'
$ws.Range("E12").Value = '  +3.38%  '
$ws.Range("D13").Value = '3.654.40'
$ws.Range("E13").Value = '  +1.50%  '
$ws.Range("D14").Value = '''25.72'
$ws.Range("E14").Value = '  +1.97%  '
$ws.Range("E15").Value = '  +3.68%  '
$ws.Range("D16").Value = '57.918.60'
$ws.Range("E16").Value = '  +1.79%  '
$ws.Range("D17").Value = '3.107.60'
$ws.Range("E17").Value = '  +1.06%  '
$ws.Range("D18").Value = '''6.15'
$ws.Range("E18").Value = '  +5.01%  '
$ws.Range("E19").Value = '  +3.38%  '
$ws.Range("D20").Value = '''8.12'
$ws.Range("E20").Value = '  +3.88%  '
$ws.Range("D21").Value = '''373.32'
$ws.Range("E21").Value = '  +7.95%  '
$ws.Range("E22").Value = '  +0.20%  '
$ws.Range("D23").Value = '''5.73'
$ws.Range("E23").Value = '  -1.08%  '
$ws.Range("D24").Value = '''69.38'
$ws.Range("E24").Value = '  +1.75%  '
$ws.Range("E25").Value = '  +2.44%  '
$ws.Range("E26").Value = '  +0.29%  '
$ws.Range("E27").Value = '  +0.18%  '
$ws.Range("D28").Value = '0.0₃0884'
$ws.Range("E28").Value = '  +3.19%  '
$ws.Range("D29").Value = '''7.59'
$ws.Range("E29").Value = '  +4.71%  '
$ws.Range("E30").Value = '  +4.97%  '
$ws.Range("E31").Value = '  +0.15%  '
$ws.Range("D32").Value = '''21.51'
$ws.Range("E32").Value = '  +3.89%  '
$ws.Range("D33").Value = '''5.15'
$ws.Range("E33").Value = '  +5.41%  '
$ws.Range("D34").Value = '''1.17'
$ws.Range("E34").Value = '  +3.53%  '
$ws.Range("D35").Value = '''160.44'
$ws.Range("E35").Value = '  +0.99%  '
$ws.Range("E36").Value = '  +3.40%  '
$ws.Range("D37").Value = '''1.31'
$ws.Range("E37").Value = '  +6.28%  '
$ws.Range("D38").Value = '''25.56'
$ws.Range("E38").Value = '  -0.10%  '
$ws.Range("E39").Value = '  +4.36%  '
$ws.Range("D40").Value = '''0.0671'
$ws.Range("E40").Value = '  +3.18%  '
$ws.Range("D41").Value = '2.558.62'
$ws.Range("E41").Value = '  +7.49%  '
$ws.Range("D42").Value = '''4.15'
$ws.Range("E42").Value = '  +3.96%  '
$ws.Range("D43").Value = '''38.40'
$ws.Range("E43").Value = '  +5.08%  '
$ws.Range("D44").Value = '''0.697'
$ws.Range("E44").Value = '  +1.28%  '
$ws.Range("D45").Value = '''0.0270'
$ws.Range("E45").Value = '  +3.13%  '
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("D47").Value = '''0.981'
$ws.Range("E47").Value = '  +2.95%  '
$ws.Range("D48").Value = '''6.14'
$ws.Range("E48").Value = '  +3.73%  '
$ws.Range("D49").Value = '''19.93'
$ws.Range("E49").Value = '  +1.88%  '
$ws.Range("D50").Value = '''0.0948'
$ws.Range("E50").Value = '  +6.79%  '
$ws.Range("D51").Value = '''0.747'
$ws.Range("E51").Value = '  -0.17%  '
